# Changes to the datasheet for check with jenkins and integration of mail functionality
#
# The "DataSheet" sheet had a long list of customer ids (rows 19-76) under the
# "deleteCustomer" / "ID" header block. That list is trimmed down to a
# shorter, refreshed set of customer ids (rows 19-38), row 39 is left blank,
# and a brand-new small "getCustomerDetails" / "ID" / <customer id> block is
# added right after it (rows 40-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

# Drop the old tail of the customer-id list (rows 39 through 76) entirely so
# the sheet's used range shrinks back down.
$ws.Range("A39:A76").EntireRow.Delete()

# Refreshed customer ids for the "deleteCustomer" block (rows 19-38).
$customerIds = @(
    "cus_Kv3NNXOk2nWmh3",
    "cus_Kv3NM8lN8vTSRd",
    "cus_Kv3NKvf9UkFTwJ",
    "cus_Kv3NsPc8mUAUHg",
    "cus_Kv3NLaOxqFK95z",
    "cus_Kv3NGBoZg9yZcZ",
    "cus_Kv3N22I2BOFbuv",
    "cus_Kv3NLfgvNOINlb",
    "cus_Kv3NfUhd7XfnxZ",
    "cus_Kv3N6tQI3FzQWF",
    "cus_Kv3NPnn1hm5Jx8",
    "cus_Kv3NDeLsrnPa8k",
    "cus_Kv3NOirrfpE4VI",
    "cus_Kv3NMP3KqwNl0k",
    "cus_Kv3NKmz1AtGURh",
    "cus_Kv3NoziMXiNZll",
    "cus_Kv3N2rFCyLLjWj",
    "cus_Kv3Njxv0unpFTw",
    "cus_Kv3Np43XgxP5w8",
    "cus_Kv3NWrwrreXumo"
)

for ($i = 0; $i -lt $customerIds.Length; $i++) {
    $row = 19 + $i
    $ws.Cells.Item($row, 1).Value = $customerIds[$i]
}

# Row 39 stays blank (gap), then a new small API section is appended.
$ws.Range("A40").Value = "getCustomerDetails"
$ws.Range("A41").Value = "ID"
$ws.Range("A42").Value = "cus_Kv3kMI4KlpKN94"

# Leave the cursor where the author left it.
$ws.Activate()
$ws.Range("A19").Select()
